$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Big "This study ..." paragraph (section "1. Advances made by this
#    study"): rewrite the whole paragraph text in one shot so Word collapses
#    all the old tiny runs into a single run, just like the target diff
#    shows. Two actual wording changes are folded in along the way:
#      a) "...the field of complex traits, primarily because it has
#          important implications for disease prediction and evolutionary
#          theory." -> "...the field of complex traits in humans and in
#          other species. It has important implications for evolutionary
#          theory, for the theory of the resemblance between relatives and
#          for applications such as trait prediction in artificial
#          selection program and for human disease."
#      b) "(as do we in the manuscript!)" -> "(as do we in the manuscript)"
# ---------------------------------------------------------------------------

$oldPara = "This study uses advanced computational methodology to answer the previously unresolved question of whether epistasis arises in human complex traits. As Reviewer 1 points out, the topic of epistasis is of considerable significance to the field of complex traits, primarily because it has important implications for disease prediction and evolutionary theory. Despite numerous highly cited reviews debating its importance over the past decade, statistical and computational limitations have made it impossible to search for epistasis empirically. Reviewer 2 correctly states that epistasis has already been shown (as do we in the manuscript!), but the crucial difference is that when epistasis has been previously reported it has been in model organisms (e.g. yeast, chickens, mice) where genetic variation is generated artificially (e.g. through gene knockout studies or extreme selective breeding). This is largely irrelevant to the topic of our study, and the presentation of credible empirical evidence for epistasis influencing human complex traits is conspicuously absent from the literature."

$newPara = "This study uses advanced computational methodology to answer the previously unresolved question of whether epistasis arises in human complex traits. As Reviewer 1 points out, the topic of epistasis is of considerable significance to the field of complex traits in humans and in other species. It has important implications for evolutionary theory, for the theory of the resemblance between relatives and for applications such as trait prediction in artificial selection program and for human disease. Despite numerous highly cited reviews debating its importance over the past decade, statistical and computational limitations have made it impossible to search for epistasis empirically. Reviewer 2 correctly states that epistasis has already been shown (as do we in the manuscript), but the crucial difference is that when epistasis has been previously reported it has been in model organisms (e.g. yeast, chickens, mice) where genetic variation is generated artificially (e.g. through gene knockout studies or extreme selective breeding). This is largely irrelevant to the topic of our study, and the presentation of credible empirical evidence for epistasis influencing human complex traits is conspicuously absent from the literature."

$d.Content.Find.Execute($oldPara, $true, $false, $false, $false, $false, $true, 1, $false, $newPara, 2)

# ---------------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the signature line to the point right
#    after "(e.g. through " in the big paragraph above -- i.e. where the
#    edit above was actually last being made, which is where Word leaves
#    its "go back to last edit" marker.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$anchor = $d.Content
$anchor.Find.Execute("generated artificially (e.g. through ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $d.Range($anchor.End, $anchor.End)
$d.Bookmarks.Add("_GoBack", $target)

Write-Output "done"
